# New crime data collected — weekly CompStat update (63rd Precinct)
# Updates: report header (volume number + week-covering dates) and the
# Week to Date / 28 Day / Year to Date / 2 Year crime-complaint figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------
# Header text: "Volume 30   Number  18" -> "...19"
#              "Report Covering the Week  5/1/2023  Through  5/7/2023"
#                -> "...5/8/2023  Through  5/14/2023"
# ---------------------------------------------------------------
$ws.Range("A8").Value = "Volume 30   Number  19"
$ws.Range("C9").Value = "Report Covering the Week  5/8/2023  Through  5/14/2023"

# ---------------------------------------------------------------
# Donor cells used to copy cell formats (number style / text style)
# when a cell's underlying type needs to flip between numeric and
# "N/A" text ("0" / "***.*") while keeping the sheet's existing
# style catalogue (no new styles created).
# ---------------------------------------------------------------
$donorTextZero = $ws.Range("C14")   # style used for text "0"
$donorTextNA   = $ws.Range("E14")   # style used for text "***.*"
$donorNumber   = $ws.Range("I14")   # plain integer style
$donorPercent  = $ws.Range("L14")   # percent-change number style

function Set-TextCell($ref, $text) {
    $dst = $ws.Range($ref)
    $dst.NumberFormat = "@"
    $dst.Value = $text
    $donorTextCell = $ws.Range($(if ($text -eq "0") { "C14" } else { "E14" }))
    $donorTextCell.Copy()
    $dst.PasteSpecial(-4122)
}

function Set-NumberCell($ref, $value, $percent) {
    $dst = $ws.Range($ref)
    $dst.Value = $value
    $donorCell = $ws.Range($(if ($percent) { "L14" } else { "I14" }))
    $donorCell.Copy()
    $dst.PasteSpecial(-4122)
}

# ---------------------------------------------------------------
# Row 14 (Murder)
# ---------------------------------------------------------------
$ws.Range("N14").Value = -77.777777777777

# ---------------------------------------------------------------
# Row 16 (Robbery)
# ---------------------------------------------------------------
Set-TextCell "D16" "0"
Set-TextCell "E16" "***.*"
$ws.Range("G16").Value = 6
$ws.Range("H16").Value = 16.666666666666
$ws.Range("I16").Value = 33
$ws.Range("K16").Value = -19.512195121951
$ws.Range("L16").Value = 73.684210526315
$ws.Range("M16").Value = -61.176470588235
$ws.Range("N16").Value = -86.307053941908

# ---------------------------------------------------------------
# Row 17 (Fel. Assault)
# ---------------------------------------------------------------
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 50
$ws.Range("F17").Value = 9
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 59
$ws.Range("J17").Value = 59
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 5.357142857142
$ws.Range("M17").Value = 31.111111111111
$ws.Range("N17").Value = -45.370370370370

# ---------------------------------------------------------------
# Row 18 (Burglary)
# ---------------------------------------------------------------
$ws.Range("C18").Value = 2
Set-NumberCell "D18" 2 $false
Set-NumberCell "E18" 0 $true
$ws.Range("F18").Value = 6
$ws.Range("G18").Value = 6
$ws.Range("I18").Value = 37
$ws.Range("J18").Value = 38
$ws.Range("K18").Value = -2.631578947368
$ws.Range("L18").Value = 15.625
$ws.Range("M18").Value = -61.855670103092
$ws.Range("N18").Value = -92.110874200426

# ---------------------------------------------------------------
# Row 19 (Gr. Larceny)
# ---------------------------------------------------------------
$ws.Range("C19").Value = 10
$ws.Range("E19").Value = -16.666666666666
$ws.Range("F19").Value = 56
$ws.Range("G19").Value = 44
$ws.Range("H19").Value = 27.272727272727
$ws.Range("I19").Value = 219
$ws.Range("J19").Value = 166
$ws.Range("K19").Value = 31.927710843373
$ws.Range("L19").Value = 99.090909090909
$ws.Range("M19").Value = 22.346368715083
$ws.Range("N19").Value = -15.116279069767

# ---------------------------------------------------------------
# Row 20 (G.L.A.)
# ---------------------------------------------------------------
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 5
$ws.Range("H20").Value = -37.5
$ws.Range("I20").Value = 38
$ws.Range("J20").Value = 40
$ws.Range("K20").Value = -5
$ws.Range("L20").Value = 26.666666666666
$ws.Range("M20").Value = -43.283582089552
$ws.Range("N20").Value = -96.451914098972

# ---------------------------------------------------------------
# Row 21 (TOTAL)
# ---------------------------------------------------------------
$ws.Range("C21").Value = 19
$ws.Range("D21").Value = 18
$ws.Range("E21").Value = 5.555555555555
$ws.Range("F21").Value = 84
$ws.Range("G21").Value = 73
$ws.Range("H21").Value = 15.068493150684
$ws.Range("I21").Value = 393
$ws.Range("J21").Value = 347
$ws.Range("K21").Value = 13.256484149855
$ws.Range("L21").Value = 57.2
$ws.Range("M21").Value = -17.782426778242
$ws.Range("N21").Value = -81.847575057736

# ---------------------------------------------------------------
# Row 23 (Housing)
# ---------------------------------------------------------------
$ws.Range("M23").Value = -60

# ---------------------------------------------------------------
# Row 24 (Petit Larceny)
# ---------------------------------------------------------------
$ws.Range("C24").Value = 24
$ws.Range("D24").Value = 23
$ws.Range("E24").Value = 4.347826086956
$ws.Range("F24").Value = 117
$ws.Range("G24").Value = 96
$ws.Range("H24").Value = 21.875
$ws.Range("I24").Value = 448
$ws.Range("J24").Value = 392
$ws.Range("K24").Value = 14.285714285714
$ws.Range("L24").Value = 62.318840579710
$ws.Range("M24").Value = 37.003058103975

# ---------------------------------------------------------------
# Row 25 (Misd. Assault)
# ---------------------------------------------------------------
$ws.Range("C25").Value = 4
$ws.Range("E25").Value = 33.333333333333
$ws.Range("F25").Value = 18
$ws.Range("H25").Value = 80
$ws.Range("I25").Value = 94
$ws.Range("J25").Value = 75
$ws.Range("K25").Value = 25.333333333333
$ws.Range("L25").Value = 9.302325581395
$ws.Range("M25").Value = -16.071428571428

# ---------------------------------------------------------------
# Row 26 (UCR Rape*)
# ---------------------------------------------------------------
$ws.Range("J26").Value = 6
$ws.Range("K26").Value = 83.333333333333

# ---------------------------------------------------------------
# Row 27 (Other Sex Crimes)
# ---------------------------------------------------------------
Set-TextCell "D27" "0"
Set-TextCell "E27" "***.*"
$ws.Range("I27").Value = 15
$ws.Range("K27").Value = 50
$ws.Range("L27").Value = 150

# ---------------------------------------------------------------
# Row 28 (Shooting Vic.)
# ---------------------------------------------------------------
Set-TextCell "F28" "0"
$ws.Range("H28").Value = -100
$ws.Range("N28").Value = -75

# ---------------------------------------------------------------
# Row 29 (Shooting Inc.)
# ---------------------------------------------------------------
Set-TextCell "F29" "0"
$ws.Range("H29").Value = -100
$ws.Range("N29").Value = -77.777777777777

# ---------------------------------------------------------------
# Row 30 (Hate Crimes)
# ---------------------------------------------------------------
Set-TextCell "F30" "0"

Write-Host "Weekly crime data refreshed."
